$d = $word.ActiveDocument
$rng = $d.Content
$found = $rng.Find.Execute("2019.ERZ.4908")
if (-not $found) {
    Write-Output "ERROR: text not found"
}
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' + `
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>' + `
'<w:p w14:paraId="4490228C" w14:textId="77777777" w:rsidR="00C476A1" w:rsidRPr="005D2900" w:rsidRDefault="002645A8" w:rsidP="005D2900"><w:pPr><w:pStyle w:val="Text85pt"/></w:pPr><w:sdt><w:sdtPr><w:tag w:val="Reference_Label"/><w:id w:val="-1816484056"/><w:placeholder><w:docPart w:val="905FF8AEF24C438183EBDE81B23DC656"/></w:placeholder><w:dataBinding w:prefixMappings="xmlns:ns=''http://schemas.officeatwork.com/CustomXMLPart''" w:xpath="/ns:officeatwork/ns:Reference_Label" w:storeItemID="{C9EF7656-0210-462C-829B-A9AFE99E1459}"/><w:text w:multiLine="1"/></w:sdtPr><w:sdtEndPr/><w:sdtContent><w:r w:rsidR="00C476A1"><w:t xml:space="preserve">Notre </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidR="00C476A1"><w:t>réf:</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidR="00C476A1"><w:t xml:space="preserve"> </w:t></w:r></w:sdtContent></w:sdt><w:r w:rsidR="00C476A1" w:rsidRPr="005D2900"><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>fallNummer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}</w:t></w:r></w:p>' + `
'</w:body></w:document>' + `
'</pkg:xmlData></pkg:part></pkg:package>'
$rng.InsertXML($xml)
Write-Output "Paragraph replaced successfully"
